$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, Price (D), Volume(1h) (E) updates per the cryptos list refresh
# Price values that look like plain numbers get a leading apostrophe so
# Excel keeps them as text (matching the sheet convention), same as the
# multi-dot price strings which are already non-numeric text.
$updates = @(
    @{ Row = 2; D = "28.165.55"; E = "  -3.18%  " }
    @{ Row = 3; D = "1.926.94"; E = "  -2.21%  " }
    @{ Row = 4; D = $null; E = "  -0.61%  " }
    @{ Row = 5; D = "'329.43"; E = "  +0.16%  " }
    @{ Row = 6; D = $null; E = "  -0.51%  " }
    @{ Row = 7; D = "'0.4729"; E = "  -5.04%  " }
    @{ Row = 8; D = "'0.4060"; E = "  -3.76%  " }
    @{ Row = 9; D = $null; E = "  -0.18%  " }
    @{ Row = 10; D = "'0.08438"; E = "  -8.87%  " }
    @{ Row = 11; D = "'1.047"; E = "  -4.86%  " }
    @{ Row = 12; D = "'22.26"; E = "  -2.68%  " }
    @{ Row = 13; D = "1.914.24"; E = "  -3.60%  " }
    @{ Row = 14; D = "'7.520"; E = "  -4.92%  " }
    @{ Row = 15; D = "'6.095"; E = "  -5.82%  " }
    @{ Row = 16; D = "'1.005"; E = "  -0.45%  " }
    @{ Row = 17; D = "'90.45"; E = "  -1.61%  " }
    @{ Row = 18; D = "'0.00001064"; E = "  -3.73%  " }
    @{ Row = 19; D = "'0.06587"; E = "  -1.94%  " }
    @{ Row = 20; D = "'18.18"; E = "  -5.75%  " }
    @{ Row = 21; D = "'1.003"; E = "  -0.59%  " }
    @{ Row = 22; D = "'5.774"; E = "  -3.23%  " }
    @{ Row = 23; D = "28.175.85"; E = "  -3.22%  " }
    @{ Row = 24; D = "'11.43"; E = "  -4.46%  " }
    @{ Row = 25; D = "'2.275"; E = "  +0.25%  " }
    @{ Row = 26; D = "2.147.54"; E = "  -3.27%  " }
    @{ Row = 27; D = "'154.67"; E = "  -0.51%  " }
    @{ Row = 28; D = "'20.10"; E = "  -2.92%  " }
    @{ Row = 29; D = "'2.153"; E = "  -4.75%  " }
    @{ Row = 30; D = "'5.750"; E = "  -9.25%  " }
    @{ Row = 31; D = "'123.80"; E = "  -2.34%  " }
    @{ Row = 32; D = "'0.9791"; E = "  -6.55%  " }
    @{ Row = 33; D = "'0.09604"; E = "  -2.54%  " }
    @{ Row = 34; D = "'1.443"; E = "  -5.06%  " }
    @{ Row = 35; D = "'5.561"; E = "  -4.61%  " }
    @{ Row = 36; D = $null; E = "  -2.31%  " }
    @{ Row = 37; D = "'9.103"; E = "  +0.46%  " }
    @{ Row = 38; D = "'0.02317"; E = "  -4.87%  " }
    @{ Row = 39; D = "'0.06177"; E = "  -3.45%  " }
    @{ Row = 40; D = "'1.241"; E = "  -6.19%  " }
    @{ Row = 41; D = "'0.6175"; E = "  -4.73%  " }
    @{ Row = 42; D = "'11.06"; E = "  -3.38%  " }
    @{ Row = 43; D = "'1.003"; E = "  -0.51%  " }
    @{ Row = 44; D = "'0.1902"; E = "  -4.50%  " }
    @{ Row = 45; D = "'1.312"; E = "  -3.14%  " }
    @{ Row = 46; D = "'0.5894"; E = "  -5.31%  " }
    @{ Row = 47; D = $null; E = "  -3.31%  " }
    @{ Row = 48; D = "'2.035"; E = "  -7.56%  " }
    @{ Row = 49; D = $null; E = "  -0.27%  " }
    @{ Row = 50; D = "'0.06809"; E = "  -2.34%  " }
    @{ Row = 51; D = "'109.91"; E = "  -2.88%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
